# Generate Report for Handback
# The "fb95c18c-37b4-4151-91d2-cf4ff29d38f6" entry has been handed back and is
# now in sync with en-US. Update the Overview status columns plus the
# per-language "Status" and "Latest Handback DateTime" fields for that row.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet: columns zh-cn (B) and de-de (C) for the fb95c18c row (row 3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# --- zh-cn sheet: Status (B3) and Latest Handback DateTime (G3) for the fb95c18c row (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $status
$wsZhCn.Range("G3").Value = "2016-03-09 10:43:28"

# --- de-de sheet: Status (B3) and Latest Handback DateTime (G3) for the fb95c18c row (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $status
$wsDeDe.Range("G3").Value = "2016-03-09 10:43:37"
